# Update cryptos list sheet with latest scraped price/volume data.
# Mirrors the upstream GitHub Actions job that refreshes cryptos.xlsx.
#
# Many Price/Volume cells hold plain text (e.g. "217.36" or "7.49") that
# Excel would otherwise auto-convert to a number when assigned directly.
# For values that parse as a number we briefly force a text number format,
# assign the literal string, then clear the temporary format again so the
# cell keeps its original (default) style while still storing text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.260.43'
$ws.Range("E2").Value = '  +1.64%  '
$ws.Range("D3").Value = '1.644.95'
$ws.Range("E3").Value = '  +0.35%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.36'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.21%  '
$ws.Range("E6").Value = '  +0.81%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  +1.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.07'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +0.84%  '
$ws.Range("E11").Value = '  +0.53%  '
$ws.Range("D12").Value = '1.875.43'
$ws.Range("E12").Value = '  +0.48%  '
$ws.Range("D13").Value = '1.639.98'
$ws.Range("E13").Value = '  -0.11%  '
$ws.Range("E14").Value = '  +0.95%  '
$ws.Range("E15").Value = '  +3.16%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '67.09'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.84%  '
$ws.Range("D17").Value = '27.256.82'
$ws.Range("D18").Value = '0.0₃0741'
$ws.Range("E18").Value = '  +1.64%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '220.32'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.51%  '
$ws.Range("E20").Value = '  +0.12%  '
$ws.Range("E21").Value = '  +3.94%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.53'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +3.82%  '
$ws.Range("E23").Value = '  +0.60%  '
$ws.Range("E24").Value = '  +0.04%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '148.80'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.01%  '
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.49'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +1.75%  '
$ws.Range("B27").Value = 'BinanceUSD'
$ws.Range("C27").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.08%  '
$ws.Range("E28").Value = '  -0.66%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.75'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.10%  '
$ws.Range("E30").Value = '  +1.87%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.18'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.10%  '
$ws.Range("E32").Value = '  +2.10%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.01'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.45%  '
$ws.Range("D34").Value = '1.308.10'
$ws.Range("E34").Value = '  +3.95%  '
$ws.Range("E35").Value = '  +1.52%  '
$ws.Range("E36").Value = '  +1.40%  '
$ws.Range("E37").Value = '  -0.10%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.554'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +3.78%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.862'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +3.49%  '
$ws.Range("E40").Value = '  +0.09%  '
$ws.Range("E41").Value = '  +0.57%  '
$ws.Range("E42").Value = '  +5.86%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.33'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.53%  '
$ws.Range("D44").Value = '1.785.10'
$ws.Range("E44").Value = '  +0.42%  '
$ws.Range("E45").Value = '  +0.55%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '91.97'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.07%  '
$ws.Range("E47").Value = '  +2.04%  '
$ws.Range("D48").Value = '0.0₆0107'
$ws.Range("E48").Value = '  +1.86%  '
$ws.Range("E49").Value = '  +0.06%  '
$ws.Range("E50").Value = '  +0.67%  '
$ws.Range("E51").Value = '  +0.45%  '
